$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 154-163: drop the placeholder empty D/E cells (dimension stays, values just vanish) ---
$cellsToClear = @("D154","E154","E155","E156","E157","E158","E159","E160","E161","E162","D163","E163")
foreach ($ref in $cellsToClear) {
    $ws.Range($ref).ClearContents() | Out-Null
}

# --- Append new log rows 164-193 (two more script runs appended to the log) ---
$rows = @(
    @{ Row=164; A="2026-02-26 17:28:29"; B="2026-02-26 18:32:19"; C="script_start"; D=$null; E=$null; F="2026-02-26 17:28:34" }
    @{ Row=165; A="2026-02-26 17:28:29"; B="2026-02-26 18:32:19"; C="split_cycle_start"; D=1; E=$null; F="2026-02-26 17:28:34" }
    @{ Row=166; A="2026-02-26 17:28:29"; B="2026-02-26 18:32:19"; C="feedallgachas_start"; D=1; E=$null; F="2026-02-26 17:28:34" }
    @{ Row=167; A="2026-02-26 17:28:29"; B="2026-02-26 18:32:19"; C="feedallgachas_end"; D=1; E=$null; F="2026-02-26 17:54:05" }
    @{ Row=168; A="2026-02-26 17:28:29"; B="2026-02-26 18:32:19"; C="collectcrack_1_6_start"; D=1; E=$null; F="2026-02-26 17:54:07" }
    @{ Row=169; A="2026-02-26 17:28:29"; B="2026-02-26 18:32:19"; C="collectcrack_1_6_end"; D=1; E=$null; F="2026-02-26 18:07:31" }
    @{ Row=170; A="2026-02-26 17:28:29"; B="2026-02-26 18:32:19"; C="feedallgachas_between_start"; D=1; E=$null; F="2026-02-26 18:07:31" }
    @{ Row=171; A="2026-02-26 17:28:29"; B="2026-02-26 18:32:19"; C="feedallgachas_between_end"; D=1; E=$null; F="2026-02-26 18:31:21" }
    @{ Row=172; A="2026-02-26 17:28:29"; B="2026-02-26 18:32:19"; C="collectcrack_7_12_start"; D=1; E=$null; F="2026-02-26 18:31:23" }
    @{ Row=173; A="2026-02-26 17:28:29"; B="2026-02-26 18:32:19"; C="script_end"; D=$null; E=$null; F="2026-02-26 18:32:19" }
    @{ Row=174; A="2026-02-26 18:59:49"; B="2026-02-26 19:14:32"; C="script_start"; D=$null; E=$null; F="2026-02-26 18:59:54" }
    @{ Row=175; A="2026-02-26 18:59:49"; B="2026-02-26 19:14:32"; C="split_cycle_start"; D=1; E=$null; F="2026-02-26 18:59:54" }
    @{ Row=176; A="2026-02-26 18:59:49"; B="2026-02-26 19:14:32"; C="feedallgachas_start"; D=1; E=$null; F="2026-02-26 18:59:54" }
    @{ Row=177; A="2026-02-26 18:59:49"; B="2026-02-26 19:14:32"; C="feedallgachas_end"; D=1; E=$null; F="2026-02-26 18:59:54" }
    @{ Row=178; A="2026-02-26 18:59:49"; B="2026-02-26 19:14:32"; C="collectcrack_1_6_start"; D=1; E=$null; F="2026-02-26 18:59:56" }
    @{ Row=179; A="2026-02-26 18:59:49"; B="2026-02-26 19:14:32"; C="collectcrack_1_6_end"; D=1; E=$null; F="2026-02-26 18:59:56" }
    @{ Row=180; A="2026-02-26 18:59:49"; B="2026-02-26 19:14:32"; C="feedallgachas_between_start"; D=1; E=$null; F="2026-02-26 18:59:56" }
    @{ Row=181; A="2026-02-26 18:59:49"; B="2026-02-26 19:14:32"; C="feedallgachas_between_end"; D=1; E=$null; F="2026-02-26 18:59:56" }
    @{ Row=182; A="2026-02-26 18:59:49"; B="2026-02-26 19:14:32"; C="collectcrack_7_12_start"; D=1; E=$null; F="2026-02-26 18:59:59" }
    @{ Row=183; A="2026-02-26 18:59:49"; B="2026-02-26 19:14:32"; C="collectcrack_7_12_end"; D=1; E=$null; F="2026-02-26 19:14:17" }
    @{ Row=184; A="2026-02-26 18:59:49"; B="2026-02-26 19:14:32"; C="split_cycle_end"; D=1; E=$null; F="2026-02-26 19:14:18" }
    @{ Row=185; A="2026-02-26 18:59:49"; B="2026-02-26 19:14:32"; C="split_cycle_start"; D=2; E=$null; F="2026-02-26 19:14:18" }
    @{ Row=186; A="2026-02-26 18:59:49"; B="2026-02-26 19:14:32"; C="feedallgachas_start"; D=2; E=$null; F="2026-02-26 19:14:18" }
    @{ Row=187; A="2026-02-26 18:59:49"; B="2026-02-26 19:14:32"; C="feedallgachas_end"; D=2; E=$null; F="2026-02-26 19:14:18" }
    @{ Row=188; A="2026-02-26 18:59:49"; B="2026-02-26 19:14:32"; C="collectcrack_1_6_start"; D=2; E=$null; F="2026-02-26 19:14:20" }
    @{ Row=189; A="2026-02-26 18:59:49"; B="2026-02-26 19:14:32"; C="collectcrack_1_6_end"; D=2; E=$null; F="2026-02-26 19:14:20" }
    @{ Row=190; A="2026-02-26 18:59:49"; B="2026-02-26 19:14:32"; C="feedallgachas_between_start"; D=2; E=$null; F="2026-02-26 19:14:20" }
    @{ Row=191; A="2026-02-26 18:59:49"; B="2026-02-26 19:14:32"; C="feedallgachas_between_end"; D=2; E=$null; F="2026-02-26 19:14:20" }
    @{ Row=192; A="2026-02-26 18:59:49"; B="2026-02-26 19:14:32"; C="collectcrack_7_12_start"; D=2; E=$null; F="2026-02-26 19:14:23" }
    @{ Row=193; A="2026-02-26 18:59:49"; B="2026-02-26 19:14:32"; C="script_end"; D=$null; E=$null; F="2026-02-26 19:14:32" }
)

foreach ($row in $rows) {
    $r = $row.Row
    $ws.Range("A$r").Value = $row.A
    $ws.Range("B$r").Value = $row.B
    $ws.Range("C$r").Value = $row.C
    if ($row.D -ne $null) { $ws.Range("D$r").Value = $row.D }
    if ($row.E -ne $null) { $ws.Range("E$r").Value = $row.E }
    $ws.Range("F$r").Value = $row.F
}

Write-Output "Log rows 154-193 updated; sheet dimension now extends to F193."
